# Timesheet updated by Gokul
# Applies the authored content changes to the "11-4-22" sheet of the
# PTW-Timesheet workbook: fills in previously-ABSENT entries, revises a
# few task/comment descriptions, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("11-4-22")

# --- Row 3 (Kirubaharan / 13-4-22 row) -------------------------------
# Was fully "ABSENT" -> now filled in with real entries.
$ws.Range("B3").Value = "Entity relationships, Entity attributes"
$ws.Range("E3").Value = "Changes in Trainee Design - 2hr,Meeting With Rafi - 1hr 30 mins, Data Models - 1hr 30 mins total hrs- 5hrs"
$ws.Range("F3").Value = "Lunch and break 1hr 15 mins, self exploration - entity frame work - 1hr total hrs - 2 hrs 15 mins"

# C3 gets new text AND a new wrap-text / left+center aligned style.
$ws.Range("C3").Value = "Completed Changes in Trainee Design, Completed changes in Trainer Design, Worked on Data Models in TMS"
$ws.Range("C3").HorizontalAlignment = -4131
$ws.Range("C3").VerticalAlignment = -4108
$ws.Range("C3").WrapText = $true

# --- Row 7 ------------------------------------------------------------
$ws.Range("C7").Value = "Reviewed acceptance criteria for Trainee ,Trainer,training coordinator,Training head, Updated MOM"
$ws.Range("E7").Value = "Discussion with team - 20 mins, Worked on entity data model for trainee and collaborated with others-120  mins,meeting with Rafi- 100 mins"
$ws.Range("F7").Value = "Prepared for code review - 1 hour , others - 90 mins"

# --- Row 10 -------------------------------------------------------------
$ws.Range("B10").Value = "Data model Relationship"
$ws.Range("C10").Value = "Entities and attributes for entire flow"
$ws.Range("D10").ClearContents()
$ws.Range("E10").Value = "Team discussion-15 mins,working with Data model-2hrs, Meeting with client (Rafi) - 1.30 hrs,Prepared data model for entire flow - 3 hrs"
$ws.Range("F10").Value = "Break - 1 hr"

# --- Row 12 -------------------------------------------------------------
$ws.Range("C12").Value = "Reviewed acceptance criteria for Trainee, Trainer, Training co ordinator, Trainnig Head. "
$ws.Range("E12").Value = "Meeting with client - 90 mins, I've decided to work on entity data model for trainee and Integrated others- 120mins, Done entity model rough - 60 mins"
$ws.Range("F12").Value = "lunch and others 130mins"

# --- Selection moves from B5 to F4 on the active sheet ------------------
$ws.Range("F4").Select() | Out-Null
